# Auto-generated edit script: applies scheduled-runner market-data updates
# to the Chocobo_Profits-style Leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 713.8
$ws.Range("I41").Value = 1000
$ws.Range("J41").Value = 682
$ws.Range("K41").Value = 1000
$ws.Range("L41").Value = 682
$ws.Range("M41").Value = -560
$ws.Range("N41").Value = -1562
$ws.Range("H53").Value = 448.1111
$ws.Range("I53").Value = 340.91666
$ws.Range("J53").Value = 533.86664
$ws.Range("K53").Value = 340.91666
$ws.Range("L53").Value = 533.86664
$ws.Range("M53").Value = 296.08334
$ws.Range("N53").Value = -1807.86664
$ws.Range("H76").Value = 3080.973
$ws.Range("I76").Value = 3057.0286
$ws.Range("K76").Value = 3057.0286
$ws.Range("M76").Value = -2742.0286
$ws.Range("H79").Value = 3080.973
$ws.Range("I79").Value = 3057.0286
$ws.Range("K79").Value = 3057.0286
$ws.Range("M79").Value = -1965.0286
$ws.Range("H98").Value = 2317.516
$ws.Range("I98").Value = 847
$ws.Range("J98").Value = 4991.1816
$ws.Range("K98").Value = 847
$ws.Range("L98").Value = 4991.1816
$ws.Range("M98").Value = 651
$ws.Range("N98").Value = -7987.1816
$ws.Range("H122").Value = 2317.516
$ws.Range("I122").Value = 847
$ws.Range("J122").Value = 4991.1816
$ws.Range("K122").Value = 2541
$ws.Range("L122").Value = 14973.5448
$ws.Range("M122").Value = -91
$ws.Range("N122").Value = -19873.5448
$ws.Range("H141").Value = 5412.4375
$ws.Range("I141").Value = 5484.023
$ws.Range("J141").Value = 4625
$ws.Range("K141").Value = 16452.069
$ws.Range("L141").Value = 13875
$ws.Range("M141").Value = -11272.069
$ws.Range("N141").Value = -24235

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2583.1538
$ws.Range("I122").Value = 1507.3636
$ws.Range("J122").Value = 8500
$ws.Range("K122").Value = 4522.0908
$ws.Range("L122").Value = 25500
$ws.Range("M122").Value = -2072.0908
$ws.Range("N122").Value = -30400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1669.5555
$ws.Range("I86").Value = 1232
$ws.Range("K86").Value = 1232
$ws.Range("M86").Value = -109
$ws.Range("H89").Value = 1669.5555
$ws.Range("I89").Value = 1232
$ws.Range("K89").Value = 6160
$ws.Range("M89").Value = -544

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1822.9867
$ws.Range("I58").Value = 1623.5
$ws.Range("J58").Value = 2774.3845
$ws.Range("K58").Value = 1623.5
$ws.Range("L58").Value = 2774.3845
$ws.Range("M58").Value = -1420.5
$ws.Range("N58").Value = -3180.3845
$ws.Range("H122").Value = 2740.5715
$ws.Range("I122").Value = 2329.7273
$ws.Range("J122").Value = 3192.5
$ws.Range("K122").Value = 6989.1819
$ws.Range("L122").Value = 9577.5
$ws.Range("M122").Value = -4539.1819
$ws.Range("N122").Value = -14477.5
$ws.Range("H136").Value = 1822.9867
$ws.Range("I136").Value = 1623.5
$ws.Range("J136").Value = 2774.3845
$ws.Range("K136").Value = 4870.5
$ws.Range("L136").Value = 8323.1535
$ws.Range("M136").Value = -2320.5
$ws.Range("N136").Value = -13423.1535

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 4333.3335
$ws.Range("J59").Value = 4333.3335
$ws.Range("L59").Value = 13000.0005
$ws.Range("N59").Value = -14080.0005
$ws.Range("H64").Value = 1771.4286
$ws.Range("I64").Value = 1000
$ws.Range("J64").Value = 1900
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 5700
$ws.Range("M64").Value = -2730
$ws.Range("N64").Value = -6240
$ws.Range("H67").Value = 1771.4286
$ws.Range("I67").Value = 1000
$ws.Range("J67").Value = 1900
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 5700
$ws.Range("M67").Value = -2064
$ws.Range("N67").Value = -7572
$ws.Range("H70").Value = 2785.4285
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 2999.6667
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 8999.000100000001
$ws.Range("M70").Value = -4185
$ws.Range("N70").Value = -9629.000100000001
$ws.Range("H73").Value = 2785.4285
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 2999.6667
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 8999.000100000001
$ws.Range("M73").Value = -3408
$ws.Range("N73").Value = -11183.0001
$ws.Range("H87").Value = 2300
$ws.Range("I87").Value = 1950
$ws.Range("J87").Value = 3000
$ws.Range("K87").Value = 5850
$ws.Range("L87").Value = 9000
$ws.Range("M87").Value = -4602
$ws.Range("N87").Value = -11496
$ws.Range("H90").Value = 2300
$ws.Range("I90").Value = 1950
$ws.Range("J90").Value = 3000
$ws.Range("K90").Value = 17550
$ws.Range("L90").Value = 27000
$ws.Range("M90").Value = -11310
$ws.Range("N90").Value = -39480
$ws.Range("H131").Value = 9091800
$ws.Range("J131").Value = 952.76
$ws.Range("L131").Value = 2858.28
$ws.Range("N131").Value = -12938.28

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6256.421
$ws.Range("I70").Value = 5751.08
$ws.Range("J70").Value = 7228.231
$ws.Range("K70").Value = 5751.08
$ws.Range("L70").Value = 7228.231
$ws.Range("M70").Value = -5481.08
$ws.Range("N70").Value = -7768.231
$ws.Range("H73").Value = 6256.421
$ws.Range("I73").Value = 5751.08
$ws.Range("J73").Value = 7228.231
$ws.Range("K73").Value = 5751.08
$ws.Range("L73").Value = 7228.231
$ws.Range("M73").Value = -4815.08
$ws.Range("N73").Value = -9100.231
$ws.Range("H122").Value = 3551.5833
$ws.Range("I122").Value = 2301.4443
$ws.Range("K122").Value = 6904.3329
$ws.Range("M122").Value = -4454.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H82").Value = 1406.4615
$ws.Range("I82").Value = 630.6667
$ws.Range("J82").Value = 2071.4285
$ws.Range("K82").Value = 630.6667
$ws.Range("L82").Value = 2071.4285
$ws.Range("M82").Value = -269.6667
$ws.Range("N82").Value = -2793.4285
$ws.Range("H85").Value = 1406.4615
$ws.Range("I85").Value = 630.6667
$ws.Range("J85").Value = 2071.4285
$ws.Range("K85").Value = 630.6667
$ws.Range("L85").Value = 2071.4285
$ws.Range("M85").Value = 617.3333
$ws.Range("N85").Value = -4567.4285
$ws.Range("H122").Value = 4122.8335
$ws.Range("I122").Value = 2714.6667
$ws.Range("J122").Value = 5531
$ws.Range("K122").Value = 8144.000100000001
$ws.Range("L122").Value = 16593
$ws.Range("M122").Value = -5694.000100000001
$ws.Range("N122").Value = -21493

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2865.875
$ws.Range("I122").Value = 1699.7778
$ws.Range("K122").Value = 5099.3334
$ws.Range("M122").Value = -2649.3334
